$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2023-12-08 Friday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2023-12-09 Saturday", 2)

# Update the division problems in the table, cell by cell to disambiguate
# duplicate values (e.g. "62÷9=" appears twice but maps to different results).
$t = $d.Tables.Item(1)

function Set-CellText($table, $row, $col, $newText) {
    $cell = $table.Cell($row, $col)
    $rng = $cell.Range
    # Trim trailing cell-marker/paragraph-end characters by shrinking the end by 1
    $rng.MoveEnd(1, -1) | Out-Null
    $rng.Text = $newText
}

# Row 1 (table row 1)
Set-CellText $t 1 1 "95÷4="
Set-CellText $t 1 2 "43÷5="
Set-CellText $t 1 3 "15÷5="
Set-CellText $t 1 4 "38÷5="
Set-CellText $t 1 5 "27÷9="

# Row 2 (table row 5)
Set-CellText $t 5 1 "92÷6="
Set-CellText $t 5 2 "69÷5="
Set-CellText $t 5 3 "91÷2="
Set-CellText $t 5 4 "47÷8="
Set-CellText $t 5 5 "15÷8="

# Row 3 (table row 9)
Set-CellText $t 9 1 "91÷6="
Set-CellText $t 9 2 "33÷8="
Set-CellText $t 9 3 "21÷7="
Set-CellText $t 9 4 "82÷7="
Set-CellText $t 9 5 "15÷8="

# Row 4 (table row 13)
Set-CellText $t 13 1 "81÷7="
Set-CellText $t 13 2 "30÷2="
Set-CellText $t 13 3 "67÷3="
Set-CellText $t 13 4 "40÷5="
Set-CellText $t 13 5 "21÷9="

# Row 5 (table row 17)
Set-CellText $t 17 1 "21÷7="
Set-CellText $t 17 2 "98÷3="
Set-CellText $t 17 3 "54÷4="
Set-CellText $t 17 4 "78÷5="
Set-CellText $t 17 5 "33÷8="
